$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update bell_start / bell_end values
$ws.Range("C2").Value = 12.5

$ws.Range("C3").Value = 16.899999999999999
$ws.Range("D3").Value = 17.8

$ws.Range("C5").Value = 1216.44
$ws.Range("D5").Value = 1217.5

$ws.Range("C6").Value = 4996.1000000000004
$ws.Range("D6").Value = 4996.6000000000004

# Select D6 as the active cell to match the saved selection state
[void]$ws.Range("D6").Select()
